$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '63.342.24'
$ws.Cells.Item(2, 5).Value = '  -1.17%  '

$ws.Cells.Item(3, 4).Value = '2.717.12'
$ws.Cells.Item(3, 5).Value = '  -1.51%  '

$ws.Cells.Item(4, 5).Value = '  +0.08%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '559.84'
$ws.Cells.Item(5, 5).Value = '  -3.00%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '156.96'
$ws.Cells.Item(6, 5).Value = '  -1.25%  '

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.13%  '

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = '0.591'
$ws.Cells.Item(8, 5).Value = '  -2.12%  '

$ws.Cells.Item(9, 5).Value = '  -2.86%  '

$ws.Cells.Item(10, 5).Value = '  -0.66%  '

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = '5.61'
$ws.Cells.Item(11, 5).Value = '  -1.94%  '

$ws.Cells.Item(12, 5).Value = '  -4.16%  '

$ws.Cells.Item(13, 4).Value = '3.198.92'
$ws.Cells.Item(13, 5).Value = '  -1.51%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '26.43'
$ws.Cells.Item(14, 5).Value = '  -2.02%  '

$ws.Cells.Item(15, 4).Value = '63.228.60'
$ws.Cells.Item(15, 5).Value = '  -0.73%  '

$ws.Cells.Item(16, 5).Value = '  -3.30%  '

$ws.Cells.Item(17, 4).Value = '2.717.88'
$ws.Cells.Item(17, 5).Value = '  -1.59%  '

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '12.17'
$ws.Cells.Item(18, 5).Value = '  +0.02%  '

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '4.66'
$ws.Cells.Item(19, 5).Value = '  -4.08%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '350.99'
$ws.Cells.Item(20, 5).Value = '  -1.92%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '6.46'
$ws.Cells.Item(21, 5).Value = '  -4.56%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.998'
$ws.Cells.Item(22, 5).Value = '  -0.11%  '

$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '0.512'
$ws.Cells.Item(23, 5).Value = '  -4.30%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '64.29'
$ws.Cells.Item(24, 5).Value = '  -1.97%  '

$ws.Cells.Item(25, 5).Value = '  -1.23%  '

$ws.Cells.Item(26, 5).Value = '  +0.12%  '

$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '8.19'
$ws.Cells.Item(27, 5).Value = '  -4.44%  '

$ws.Cells.Item(28, 4).Value = '0.0₃0888'
$ws.Cells.Item(28, 5).Value = '  -2.18%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '1.37'
$ws.Cells.Item(29, 5).Value = '  +9.83%  '

$ws.Cells.Item(30, 5).Value = '  -0.95%  '

$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = '7.14'
$ws.Cells.Item(31, 5).Value = '  -2.59%  '

$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '166.20'
$ws.Cells.Item(32, 5).Value = '  -1.77%  '

$ws.Cells.Item(33, 5).Value = '  -0.52%  '

$ws.Cells.Item(34, 2).Value = 'USDe'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.998'
$ws.Cells.Item(34, 5).Value = '  -0.02%  '

$ws.Cells.Item(35, 2).Value = 'EthereumClassic'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '19.82'
$ws.Cells.Item(35, 5).Value = '  -2.19%  '

$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = '4.82'
$ws.Cells.Item(36, 5).Value = '  -2.57%  '

$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = '1.77'
$ws.Cells.Item(37, 5).Value = '  -2.46%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '347.04'
$ws.Cells.Item(38, 5).Value = '  -0.21%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '0.954'
$ws.Cells.Item(39, 5).Value = '  -5.16%  '

$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = '6.07'
$ws.Cells.Item(40, 5).Value = '  -4.97%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '4.05'
$ws.Cells.Item(41, 5).Value = '  -3.54%  '

$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '38.34'
$ws.Cells.Item(42, 5).Value = '  -2.12%  '

$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = '21.38'
$ws.Cells.Item(43, 5).Value = '  -2.13%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '20.77'
$ws.Cells.Item(44, 5).Value = '  -3.48%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.0573'
$ws.Cells.Item(45, 5).Value = '  -3.10%  '

$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = '0.626'
$ws.Cells.Item(46, 5).Value = '  -1.18%  '

$ws.Cells.Item(47, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '0.998'
$ws.Cells.Item(47, 5).Value = '  -0.02%  '

$ws.Cells.Item(48, 2).Value = 'Aave'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '131.36'
$ws.Cells.Item(48, 5).Value = '  -3.84%  '

$ws.Cells.Item(49, 2).Value = 'Stellar'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '0.0985'
$ws.Cells.Item(49, 5).Value = '  -3.36%  '

$ws.Cells.Item(50, 2).Value = 'VeChain'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '0.0245'
$ws.Cells.Item(50, 5).Value = '  -4.05%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '11.04'
$ws.Cells.Item(51, 5).Value = '  +0.03%  '
